# Fix link in powerpoint presentation
#
# The slide that describes the open-source Python RIFT implementation
# shows a hyperlinked text run that still points at the old repo name
# ("rift-fsm"). The repo was renamed to "rift-python", so the *visible*
# run text needs to change from:
#
#   https://github.com/brunorijsman/rift-fsm
#
# to two runs whose combined text reads:
#
#   https://github.com/brunorijsman/rift-python
#
# (split as "https://github.com/brunorijsman" + "/rift-python"), while
# keeping the existing hyperlink relationship (rId2) on both pieces.

$p = $ppt.ActivePresentation

$oldLinkText  = "https://github.com/brunorijsman/rift-fsm"
$newTextPart1 = "https://github.com/brunorijsman"
$newTextPart2 = "/rift-python"

$updated = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange

        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi)

            for ($ri = 1; $ri -le $para.Runs().Count; $ri++) {
                $run = $para.Runs($ri)

                # Paragraph/run text coming back from this object model is
                # terminated with a trailing CR (paragraph mark); strip it
                # before comparing.
                $runText = $run.Text.TrimEnd("`r")

                if ($runText -eq $oldLinkText) {
                    # Shrink the existing (hyperlinked) run to just the
                    # unchanged prefix of the URL, keeping its formatting
                    # (including the rId2 hyperlink) intact ...
                    $run.Text = $newTextPart1

                    # ... then append the new suffix as a following run.
                    # InsertAfter clones the formatting (and hyperlink) of
                    # the run it is called on, so the new "/rift-python"
                    # text stays clickable with the same link.
                    $null = $run.InsertAfter($newTextPart2)

                    $updated = $true
                }
            }
        }
    }
}

if ($updated) {
    Write-Host "Updated RIFT hackathon repo link text to rift-python."
} else {
    Write-Host "No matching link text found; presentation left unchanged."
}
